$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = $ws.Range("A2").NumberFormat()

$rows = @(
  @{ r=191; a=44249; b="0221443"; bText=$true; c=3011; d="Order 0221443 Card(Stripe)"; e=$null; f=502.68 },
  @{ r=192; a=44249; b="0221443"; bText=$true; c=2611; d="Order 0221443 Card(Stripe)"; e=$null; f=60.32 },
  @{ r=193; a=44249; b="0221443"; bText=$true; c=1930; d="Order 0221443 Card(Stripe)"; e=563; f=$null },
  @{ r=194; a=44251; b="0241219"; bText=$true; c=3011; d="Order 0241219 Swish +46723656673"; e=$null; f=691.07 },
  @{ r=195; a=44251; b="0241219"; bText=$true; c=2611; d="Order 0241219 Swish +46723656673"; e=$null; f=82.93 },
  @{ r=196; a=44251; b="0241219"; bText=$true; c=1930; d="Order 0241219 Swish +46723656673"; e=774; f=$null },
  @{ r=197; a=44252; b=6251432; bText=$false; c=3011; d="Order 6251432 Card(Stripe)"; e=$null; f=1085.71 },
  @{ r=198; a=44252; b=6251432; bText=$false; c=2611; d="Order 6251432 Card(Stripe)"; e=$null; f=130.29 },
  @{ r=199; a=44252; b=6251432; bText=$false; c=1930; d="Order 6251432 Card(Stripe)"; e=1216; f=$null },
  @{ r=200; a=44254; b=5271529; bText=$false; c=3011; d="Order 5271529 Swish +46705558420"; e=$null; f=616.07 },
  @{ r=201; a=44254; b=5271529; bText=$false; c=2611; d="Order 5271529 Swish +46705558420"; e=$null; f=73.93 },
  @{ r=202; a=44254; b=5271529; bText=$false; c=1930; d="Order 5271529 Swish +46705558420"; e=690; f=$null },
  @{ r=203; a=44254; b=9272123; bText=$false; c=3011; d="Order 9272123 Card(Stripe)"; e=$null; f=537.5 },
  @{ r=204; a=44254; b=9272123; bText=$false; c=2611; d="Order 9272123 Card(Stripe)"; e=$null; f=64.5 },
  @{ r=205; a=44254; b=9272123; bText=$false; c=1930; d="Order 9272123 Card(Stripe)"; e=602; f=$null },
  @{ r=206; a=44253; b=$null; bText=$false; c=4010; d="M&S RB BROMMA K0135"; e=2009.15; f=$null },
  @{ r=207; a=44253; b=$null; bText=$false; c=2645; d="M&S RB BROMMA K0135"; e=241.1; f=$null },
  @{ r=208; a=44253; b=$null; bText=$false; c=1930; d="M&S RB BROMMA K0135"; e=$null; f=2250.25 },
  @{ r=209; a=44253; b=$null; bText=$false; c=5460; d="IKEA BARKARBY K0135"; e=149.6; f=$null },
  @{ r=210; a=44253; b=$null; bText=$false; c=2641; d="IKEA BARKARBY K0135"; e=37.4; f=$null },
  @{ r=211; a=44253; b=$null; bText=$false; c=1930; d="IKEA BARKARBY K0135"; e=$null; f=187 },
  @{ r=212; a=44254; b=$null; bText=$false; c=4010; d="SNABBGROSS SOLNA K0135"; e=468.6; f=$null },
  @{ r=213; a=44254; b=$null; bText=$false; c=2645; d="SNABBGROSS SOLNA K0135"; e=56.23; f=$null },
  @{ r=214; a=44254; b=$null; bText=$false; c=1930; d="SNABBGROSS SOLNA K0135"; e=$null; f=524.83 },
  @{ r=215; a=44254; b=$null; bText=$false; c=5670; d="ST1 V#LLINGBY K0135"; e=785.47; f=$null },
  @{ r=216; a=44254; b=$null; bText=$false; c=2641; d="ST1 V#LLINGBY K0135"; e=196.37; f=$null },
  @{ r=217; a=44254; b=$null; bText=$false; c=1930; d="ST1 V#LLINGBY K0135"; e=$null; f=981.84 },
  @{ r=218; a=44255; b=3281219; bText=$false; c=3011; d="Order 3281219 Card(Stripe)"; e=$null; f=769.64 },
  @{ r=219; a=44255; b=3281219; bText=$false; c=2611; d="Order 3281219 Card(Stripe)"; e=$null; f=92.36 },
  @{ r=220; a=44255; b=3281219; bText=$false; c=1930; d="Order 3281219 Card(Stripe)"; e=862; f=$null },
  @{ r=221; a=44255; b=3282108; bText=$false; c=3011; d="Order 3282108 Swish +46707678891"; e=$null; f=547.32 },
  @{ r=222; a=44255; b=3282108; bText=$false; c=2611; d="Order 3282108 Swish +46707678891"; e=$null; f=65.68 },
  @{ r=223; a=44255; b=3282108; bText=$false; c=1930; d="Order 3282108 Swish +46707678891"; e=613; f=$null },
  @{ r=224; a=44255; b=$null; bText=$false; c=4010; d="WILLYS STOCKHOLM VINST K0135"; e=22.23; f=$null },
  @{ r=225; a=44255; b=$null; bText=$false; c=2645; d="WILLYS STOCKHOLM VINST K0135"; e=2.67; f=$null },
  @{ r=226; a=44255; b=$null; bText=$false; c=1930; d="WILLYS STOCKHOLM VINST K0135"; e=$null; f=24.9 },
  @{ r=227; a=44255; b=$null; bText=$false; c=4010; d="SNABBGROSS SOLNA K0135"; e=694.38; f=$null },
  @{ r=228; a=44255; b=$null; bText=$false; c=2645; d="SNABBGROSS SOLNA K0135"; e=83.33; f=$null },
  @{ r=229; a=44255; b=$null; bText=$false; c=1930; d="SNABBGROSS SOLNA K0135"; e=$null; f=777.71 },
  @{ r=230; a=44255; b=$null; bText=$false; c=4010; d="NGROCERIES AB K0135"; e=379.64; f=$null },
  @{ r=231; a=44255; b=$null; bText=$false; c=2645; d="NGROCERIES AB K0135"; e=45.56; f=$null },
  @{ r=232; a=44255; b=$null; bText=$false; c=1930; d="NGROCERIES AB K0135"; e=$null; f=425.2 }
)

foreach ($row in $rows) {
  $aCell = $ws.Cells.Item($row.r, 1)
  $aCell.Value = $row.a
  $aCell.NumberFormat = $dateFmt

  if ($row.b -ne $null) {
    $bCell = $ws.Cells.Item($row.r, 2)
    if ($row.bText) {
      $bCell.NumberFormat = "@"
    }
    $bCell.Value = $row.b
  }

  $ws.Cells.Item($row.r, 3).Value = $row.c
  $ws.Cells.Item($row.r, 4).Value = $row.d

  if ($row.e -ne $null) {
    $ws.Cells.Item($row.r, 5).Value = $row.e
  }
  if ($row.f -ne $null) {
    $ws.Cells.Item($row.r, 6).Value = $row.f
  }
}

Write-Host "Added" $rows.Count "rows"